# Mark the federal holiday columns (Veterans Day - Nov 11 on sheet "1-15",
# Thanksgiving - Nov 25 on sheet "16-End") the same way the existing
# weekend (SAT/SUN) columns are marked: narrow columns, gray-shaded cells,
# and an "X" filled in on the two shaded rows per work-site block.

$wb = $excel.ActiveWorkbook

# xlPasteFormats = -4122 ; used so the destination cells pick up the exact
# existing style entries (border/font/fill) already defined for the
# weekend columns instead of Excel inventing brand-new style records.
$xlPasteFormats = -4122

# Rows that get an "X" filled in (the two shaded data rows for each of the
# 7 repeating site blocks + the final row of the last block).
$xRows = @(5,6,8,9,11,12,14,15,17,18,20,21,23,24,26,27)

# --- Sheet "1-15": Thursday Nov 11 (Veterans Day) lives in columns V:W ---
$ws1 = $wb.Worksheets.Item("1-15")

$ws1.Columns.Item(22).ColumnWidth = 1.67   # column V -> stored width 2.5
$ws1.Columns.Item(23).ColumnWidth = 1.67   # column W -> stored width 2.5

# Column L:M (SAT) already carries the "holiday-style" shading for every
# row 2-27; copy that formatting onto V:W in one shot.
$ws1.Range("L2:M27").Copy()
$ws1.Range("V2:W27").PasteSpecial($xlPasteFormats)
$ws1.Application.CutCopyMode = $false

foreach ($r in $xRows) {
    $ws1.Range("V$r").Value = "X"
    $ws1.Range("W$r").Value = "X"
}

# --- Sheet "16-End": Thursday Nov 25 (Thanksgiving) lives in columns T:U ---
$ws2 = $wb.Worksheets.Item("16-End")

$ws2.Columns.Item(20).ColumnWidth = 1.67   # column T -> stored width 2.5
$ws2.Columns.Item(21).ColumnWidth = 1.67   # column U -> stored width 2.5

# Column J:K (SAT) carries the same holiday-style shading on this sheet.
$ws2.Range("J2:K27").Copy()
$ws2.Range("T2:U27").PasteSpecial($xlPasteFormats)
$ws2.Application.CutCopyMode = $false

foreach ($r in $xRows) {
    $ws2.Range("T$r").Value = "X"
    $ws2.Range("U$r").Value = "X"
}
